$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hydropower plant parameters")

# --- Insert two new rows at position 33 (pushes the pumped-storage rows down to 35-40) ---
$ws.Rows.Item(33).Insert()
$ws.Rows.Item(33).Insert()

# --- Row 33: year_calibration_start ---
$ws.Cells.Item(33,1).Value = "year_calibration_start"

$descStart = "[leave empty if no preference - full period will be used] the first year of the multiannual period whose discharge average is to be used for the calibration of the conventional operation"
$cellB33 = $ws.Cells.Item(33,2)
$cellB33.Value = $descStart
$r1 = $cellB33.Characters(1,57)
$r1.Font.Bold = $true
$r1.Font.Italic = $true
$r1.Font.Size = 11
$r1.Font.Name = "Calibri"
$r2 = $cellB33.Characters(58,1)
$r2.Font.Italic = $true
$r2.Font.Size = 11
$r2.Font.Name = "Calibri"
$r3 = $cellB33.Characters(59,128)
$r3.Font.Size = 11
$r3.Font.Name = "Calibri"

$ws.Rows.Item(33).RowHeight = 43.5

# --- Row 34: year_calibration_end ---
$ws.Cells.Item(34,1).Value = "year_calibration_end"

$descEnd = "[leave empty if no preference - full period will be used] the last year of the multiannual period whose discharge average is to be used for the calibration of the conventional operation"
$cellB34 = $ws.Cells.Item(34,2)
$cellB34.Value = $descEnd
$r4 = $cellB34.Characters(1,57)
$r4.Font.Bold = $true
$r4.Font.Italic = $true
$r4.Font.Size = 11
$r4.Font.Name = "Calibri"
$r5 = $cellB34.Characters(58,133)
$r5.Font.Size = 11
$r5.Font.Name = "Calibri"

$ws.Rows.Item(34).RowHeight = 43.5

# --- Row height updates caused by the wider column B (text re-wraps) ---
$ws.Rows.Item(9).RowHeight = 29
$ws.Rows.Item(22).RowHeight = 14.5
$ws.Rows.Item(25).RowHeight = 29
$ws.Rows.Item(28).RowHeight = 14.5
$ws.Rows.Item(29).RowHeight = 14.5

# --- Widen column B (engine snaps column width to a pixel grid; 62 is the
#     ColumnWidth input that lands closest to the target stored width of
#     62.81640625) ---
$ws.Columns.Item(2).ColumnWidth = 62

# --- Update the active selection shown when the file is reopened ---
$ws.Range("C7").Select()
